# "Starting the ending of a raid." — add a new column to the raids import
# sheet recording the specialty-reward item type dropped by the raid boss.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cell (column G, header row 1) and its value for the single
# data row (row 2), matching the existing name/story/raid_boss_id/... layout.
$ws.Range("G1").Value = "item_specialty_reward_type"
$ws.Range("G2").Value = "Pirate Lord Leather"

# Match the look of the other bestFit/auto-sized columns on this sheet.
$ws.Columns.Item(7).AutoFit()
